$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.963.31"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.104.17"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.87"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.82"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.103.36"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.35"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.10"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "3.622.25"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "66.929.83"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.98"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "3.107.24"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.26"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.74"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.689"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.63"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.62"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.23"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.04"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.06"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "0.0₃0939"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.75"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.943"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.15"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.29"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.67"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").Value = "2.792.25"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "371.61"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.44"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0343"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.90"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
